# repull data, push all data, mean calculation
# Update the dSF (column F) values for the weigh-in log rows with freshly
# repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 5
    4  = 2
    5  = -8
    7  = 2
    8  = 1
    9  = -1
    10 = 1
    11 = -7
    12 = -1
    13 = -5
    14 = -8
    15 = 1
    16 = -4
    17 = 1
    18 = 1
    19 = -4
    20 = 1
    22 = 4
    23 = 6
    24 = 1
    25 = 4
    26 = 2
    27 = -4
    28 = 1
    30 = 0
    32 = 3
    33 = 2
    34 = 2
    35 = -2
    37 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
